# Apply "added harvard case classification" update:
# Recompute the *_old metrics for the columns that needed the
# Harvard case classification (Ada_old=C, Avey_old=F, Buoy_old=I,
# WebMD_old=Q, doctor_TH_old=W) across rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - precision
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("I2").Value = 1
$ws.Range("Q2").Value = 0.6666666666666666
$ws.Range("W2").Value = 1

# Row 3 - recall
$ws.Range("C3").Value = 0.8
$ws.Range("F3").Value = 0.8
$ws.Range("I3").Value = 0.2
$ws.Range("Q3").Value = 0.8
$ws.Range("W3").Value = 0.6

# Row 4 - f1-score
$ws.Range("C4").Value = 0.888888888888889
$ws.Range("F4").Value = 0.7272727272727272
$ws.Range("I4").Value = 0.3333333333333334
$ws.Range("Q4").Value = 0.7272727272727272
$ws.Range("W4").Value = 0.7499999999999999

# Row 5 - f2-score
$ws.Range("C5").Value = 0.8333333333333334
$ws.Range("F5").Value = 0.7692307692307692
$ws.Range("I5").Value = 0.2380952380952381
$ws.Range("Q5").Value = 0.7692307692307692
$ws.Range("W5").Value = 0.6521739130434783

# Row 6 - NDCG
$ws.Range("C6").Value = 0.9503498753418869
$ws.Range("F6").Value = 0.6615624441563611
$ws.Range("I6").Value = 0.6791866504595295
$ws.Range("Q6").Value = 0.5157437785053824
$ws.Range("W6").Value = 0.7869036435967537

# Row 7 - M1
$ws.Range("C7").Value = $true
$ws.Range("I7").Value = $true
$ws.Range("W7").Value = $true

# Row 8 - M3
$ws.Range("C8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("I8").Value = $true
$ws.Range("W8").Value = $true

# Row 9 - M5
$ws.Range("C9").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("I9").Value = $true
$ws.Range("Q9").Value = $true
$ws.Range("W9").Value = $true

# Row 10 - position
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 3
$ws.Range("I10").Value = 1
$ws.Range("Q10").Value = 4
$ws.Range("W10").Value = 1
